# The document contains a transcribed passage with inline pseudo-tags:
#   ... tu mectras <env>au soleil bien fort</env> &amp; <env>au serain</env> l<ms>...
# The edit removes the <env>...</env> markup wrapping "au serain" while keeping
# the plain text, and tidies the spacing around the "&amp;" run:
#   ... tu mectras <env>au soleil bien fort</env> &amp; au serain</env> l<ms>...

$d = $word.ActiveDocument

# Locate the anchor text once, from the start of the document, so we find the
# one unique occurrence of this run of pseudo-tags.
$anchor = $d.Content
$found = $anchor.Find.Execute("</env> &amp; <env>au serain</env>", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not find target passage"
}
$segStart = $anchor.Start
$segEnd = $anchor.End

# Walk forward through the six runs that make up this segment, using Find
# (restricted to the segment) to recover the exact run boundaries:
#   </env> | " &" | "amp;" | " " | <env> | au serain
$pos = $segStart
$bounds = @{}
$pieces = @("</env>", " &", "amp;", " ", "<env>", "au serain")
foreach ($piece in $pieces) {
    $r = $d.Range($pos, $segEnd)
    $ok = $r.Find.Execute($piece, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        throw "Could not find piece [$piece]"
    }
    $bounds[$piece] = @($r.Start, $r.End)
    $pos = $r.End
}

# Apply the edits from the end of the segment backwards, so earlier offsets
# in $bounds stay valid (editing text shifts everything after it).

# 3) The " " run, the "<env>" run, and the "au serain" run collapse into a
#    single plain run reading " au serain" (the "<env>" tag text is deleted
#    and its neighbouring plain-text runs are merged).
$mergedStart = $bounds[" "][0]
$mergedEnd = $bounds["au serain"][1]
$merged = $d.Range($mergedStart, $mergedEnd)
$merged.Text = " au serain"

# "amp;" run (gray, Courier New) is left untouched.

# 2) " &" run (plain black) -> becomes just "&".
$r2 = $d.Range($bounds[" &"][0], $bounds[" &"][1])
$r2.Text = "&"

# 1) First "</env>" run (blue, Courier New) -> becomes a single space,
#    keeping its own run formatting.
$r1 = $d.Range($bounds["</env>"][0], $bounds["</env>"][1])
$r1.Text = " "
